# Rewrites the single opening paragraph into the expanded/edited version
# (split across many runs, with a couple of extra inserted phrases) and
# splits off a brand-new second paragraph ("Away from work ...") that
# takes over the _GoBack bookmark which used to sit inside paragraph 1.
#
# Word's Range.InsertXML (when called on a Range that spans the *entire*
# document body) replaces the whole <w:body> content while leaving the
# final <w:sectPr> alone, which lets us define the exact run/paragraph
# layout required instead of relying on Find/Replace (which would keep
# merging everything back into a single run).

$d = $word.ActiveDocument
$full = $d.Content

$newBody = @'
<w:p><w:r><w:t>Mark has joined the Customer Experience</w:t></w:r><w:r><w:t xml:space="preserve"> Buying &amp; Acquisition</w:t></w:r><w:r><w:t xml:space="preserve"> team to bring knowledge and experience in the area of test automation.  His background is strongly in software product quality through a number of years developing software</w:t></w:r><w:r><w:t xml:space="preserve"> field</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> developing, testing</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t xml:space="preserve">managing </w:t></w:r><w:r><w:t xml:space="preserve">some </w:t></w:r><w:r><w:t>off-shore</w:t></w:r><w:r><w:t xml:space="preserve"> teams.  Since returning from Telecom and project management, his focus has been in Web based products.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Away from work Mark wishes he could be bicycling or swimming several hours a day, alas there are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>chores  :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> )</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$bodyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newBody + '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$full.InsertXML($bodyXml)
